$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.85129999999999
$ws.Range("E5").Value = 12.36819999999999
$ws.Range("E9").Value = 13.01470000000001
$ws.Range("E11").Value = 13.4942
$ws.Range("B21").Value = 5.781399999999995
$ws.Range("E21").Value = 13.10089999999999
$ws.Range("B23").Value = 5.893099999999997
$ws.Range("B25").Value = 5.779299999999994
